$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new column at A, shifting Car Model..Total Sale from A-H to B-I
$ws.Range("A1").EntireColumn.Insert()

# New column A header + style (match the other header cells)
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$ws.Range("A1").Value = "months"

# Fill months for each data row
$ws.Range("A2").Value = "Januray"
$ws.Range("A3").Value = "Februrary"
$ws.Range("A4").Value = "januray"
$ws.Range("A5").Value = "March"
$ws.Range("A6").Value = "June"
$ws.Range("A7").Value = "July"
$ws.Range("A8").Value = "August"
$ws.Range("A9").Value = "januray"
$ws.Range("A10").Value = "October"
$ws.Range("A11").Value = "januray"
$ws.Range("A12").Value = "November"
$ws.Range("A13").Value = "Decenber"
$ws.Range("A14").Value = "januray"
$ws.Range("A15").Value = "March"
$ws.Range("A16").Value = "April"
$ws.Range("A17").Value = "januray"
$ws.Range("A18").Value = "Februrary"
$ws.Range("A19").Value = "March"
$ws.Range("A20").Value = "November"
$ws.Range("A21").Value = "Decenber"
$ws.Range("A22").Value = "August"
$ws.Range("A23").Value = "June"
$ws.Range("A24").Value = "September"
$ws.Range("A25").Value = "Decenber"
$ws.Range("A26").Value = "januray"
$ws.Range("A27").Value = "April"
$ws.Range("A28").Value = "September"
$ws.Range("A29").Value = "July"
$ws.Range("A30").Value = "March"
$ws.Range("A31").Value = "januray"

# After the column insert, H holds "Total Tax" (old column G) and I holds
# "Total Sale" (old column H). The target layout wants Total Sale in H and
# Total Tax in I, so swap the header labels and all data values between H and I.
$ws.Range("H1").Value = "Total Sale"
$ws.Range("I1").Value = "Total Tax"
$tmp = $ws.Range("H2").Value2
$ws.Range("H2").Value = $ws.Range("I2").Value2
$ws.Range("I2").Value = $tmp
$tmp = $ws.Range("H3").Value2
$ws.Range("H3").Value = $ws.Range("I3").Value2
$ws.Range("I3").Value = $tmp
$tmp = $ws.Range("H4").Value2
$ws.Range("H4").Value = $ws.Range("I4").Value2
$ws.Range("I4").Value = $tmp
$tmp = $ws.Range("H5").Value2
$ws.Range("H5").Value = $ws.Range("I5").Value2
$ws.Range("I5").Value = $tmp
$tmp = $ws.Range("H6").Value2
$ws.Range("H6").Value = $ws.Range("I6").Value2
$ws.Range("I6").Value = $tmp
$tmp = $ws.Range("H7").Value2
$ws.Range("H7").Value = $ws.Range("I7").Value2
$ws.Range("I7").Value = $tmp
$tmp = $ws.Range("H8").Value2
$ws.Range("H8").Value = $ws.Range("I8").Value2
$ws.Range("I8").Value = $tmp
$tmp = $ws.Range("H9").Value2
$ws.Range("H9").Value = $ws.Range("I9").Value2
$ws.Range("I9").Value = $tmp
$tmp = $ws.Range("H10").Value2
$ws.Range("H10").Value = $ws.Range("I10").Value2
$ws.Range("I10").Value = $tmp
$tmp = $ws.Range("H11").Value2
$ws.Range("H11").Value = $ws.Range("I11").Value2
$ws.Range("I11").Value = $tmp
$tmp = $ws.Range("H12").Value2
$ws.Range("H12").Value = $ws.Range("I12").Value2
$ws.Range("I12").Value = $tmp
$tmp = $ws.Range("H13").Value2
$ws.Range("H13").Value = $ws.Range("I13").Value2
$ws.Range("I13").Value = $tmp
$tmp = $ws.Range("H14").Value2
$ws.Range("H14").Value = $ws.Range("I14").Value2
$ws.Range("I14").Value = $tmp
$tmp = $ws.Range("H15").Value2
$ws.Range("H15").Value = $ws.Range("I15").Value2
$ws.Range("I15").Value = $tmp
$tmp = $ws.Range("H16").Value2
$ws.Range("H16").Value = $ws.Range("I16").Value2
$ws.Range("I16").Value = $tmp
$tmp = $ws.Range("H17").Value2
$ws.Range("H17").Value = $ws.Range("I17").Value2
$ws.Range("I17").Value = $tmp
$tmp = $ws.Range("H18").Value2
$ws.Range("H18").Value = $ws.Range("I18").Value2
$ws.Range("I18").Value = $tmp
$tmp = $ws.Range("H19").Value2
$ws.Range("H19").Value = $ws.Range("I19").Value2
$ws.Range("I19").Value = $tmp
$tmp = $ws.Range("H20").Value2
$ws.Range("H20").Value = $ws.Range("I20").Value2
$ws.Range("I20").Value = $tmp
$tmp = $ws.Range("H21").Value2
$ws.Range("H21").Value = $ws.Range("I21").Value2
$ws.Range("I21").Value = $tmp
$tmp = $ws.Range("H22").Value2
$ws.Range("H22").Value = $ws.Range("I22").Value2
$ws.Range("I22").Value = $tmp
$tmp = $ws.Range("H23").Value2
$ws.Range("H23").Value = $ws.Range("I23").Value2
$ws.Range("I23").Value = $tmp
$tmp = $ws.Range("H24").Value2
$ws.Range("H24").Value = $ws.Range("I24").Value2
$ws.Range("I24").Value = $tmp
$tmp = $ws.Range("H25").Value2
$ws.Range("H25").Value = $ws.Range("I25").Value2
$ws.Range("I25").Value = $tmp
$tmp = $ws.Range("H26").Value2
$ws.Range("H26").Value = $ws.Range("I26").Value2
$ws.Range("I26").Value = $tmp
$tmp = $ws.Range("H27").Value2
$ws.Range("H27").Value = $ws.Range("I27").Value2
$ws.Range("I27").Value = $tmp
$tmp = $ws.Range("H28").Value2
$ws.Range("H28").Value = $ws.Range("I28").Value2
$ws.Range("I28").Value = $tmp
$tmp = $ws.Range("H29").Value2
$ws.Range("H29").Value = $ws.Range("I29").Value2
$ws.Range("I29").Value = $tmp
$tmp = $ws.Range("H30").Value2
$ws.Range("H30").Value = $ws.Range("I30").Value2
$ws.Range("I30").Value = $tmp
$tmp = $ws.Range("H31").Value2
$ws.Range("H31").Value = $ws.Range("I31").Value2
$ws.Range("I31").Value = $tmp
